$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old data rows (2-15) so stale shared-string entries are dropped
# and the new values get (re)written by the user in the intended order.
$ws.Range("A2:F15").ClearContents() | Out-Null

# New data entered by the user (rows 2-9)
$data = @(
    @("Hombre", 20, "CABA", 5, "Universitario", 100000),
    @("Hombre", 21, "CABA", 3, "Universitario", 110000),
    @("Hombre", 22, "CABA", 2, "Universitario", 120000),
    @("Hombre", 20, "CABA", 5, "Universitario", 100000),
    @("Hombre", 21, "CABA", 3, "Universitario", 110000),
    @("Hombre", 22, "CABA", 2, "Universitario", 120000),
    @("Mujer",  21, "CABA", 3, "Universitario", 20000),
    @("Mujer",  22, "CABA", 2, "Universitario", 10000)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row++
}

# Widen the newly-significant columns (Years of experience / Education level)
$ws.Columns.Item(4).ColumnWidth = 17.417
$ws.Columns.Item(5).ColumnWidth = 15.584

# Update the active selection to match the final cursor location
$ws.Range("F9").Select() | Out-Null
